# Regenerate merged AHB files
# - Rename the diff-header columns from the "_old"/"_new" suffix convention
#   to the versioned "_FV2404"/"_FV2410" convention.
# - Turn the data range into a real Excel Table ("Table1").
# - Freeze the header row (top row) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old header text -> new header text for the header row (row 1).
$map = @{
    "Segmentname_old"          = "Segmentname_FV2404"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2404"
    "Segment_old"              = "Segment_FV2404"
    "Datenelement_old"         = "Datenelement_FV2404"
    "Segment ID_old"           = "Segment ID_FV2404"
    "Code_old"                 = "Code_FV2404"
    "Qualifier_old"            = "Qualifier_FV2404"
    "Beschreibung_old"         = "Beschreibung_FV2404"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2404"
    "Bedingung_old"            = "Bedingung_FV2404"
    "Segmentname_new"          = "Segmentname_FV2410"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2410"
    "Segment_new"              = "Segment_FV2410"
    "Datenelement_new"         = "Datenelement_FV2410"
    "Segment ID_new"           = "Segment ID_FV2410"
    "Code_new"                 = "Code_FV2410"
    "Qualifier_new"            = "Qualifier_FV2410"
    "Beschreibung_new"         = "Beschreibung_FV2410"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2410"
    "Bedingung_new"            = "Bedingung_FV2410"
}

$lastCol = $ws.UsedRange.Columns.Count
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}

# Convert the used range (header + 58 data rows, A1:U59) into an Excel table.
$lastRow = $ws.UsedRange.Rows.Count
$tableRange = $ws.Range("A1").Resize($lastRow, $lastCol)
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the top (header) row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
